# Add two new region rows (A208/Arapiraca 2/NE and A216/Guaruva/SC)
# to the bottom of the "REGIAO" table on Planilha1, and keep the
# table's AutoFilter range, the workbook's hidden _FilterDatabase
# defined name, and the active selection in sync with the new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new data rows -------------------------------------------------
$ws.Cells.Item(170, 1).Value = "A208"
$ws.Cells.Item(170, 2).Value = "Arapiraca 2"
$ws.Cells.Item(170, 3).Value = "NE"

$ws.Cells.Item(171, 1).Value = "A216"
$ws.Cells.Item(171, 2).Value = "Guaruva"
$ws.Cells.Item(171, 3).Value = "SC"

# --- Re-apply the AutoFilter so its range grows to cover the new rows --------
$ws.AutoFilterMode = $false
$ws.Range("A1:C171").AutoFilter()

# --- Update the hidden _FilterDatabase defined name to match -----------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Planilha1!_FilterDatabase") {
        $n.RefersTo = "=Planilha1!`$A`$1:`$C`$171"
    }
}

# --- Move the active selection to the next empty row, like after data entry --
$ws.Cells.Item(172, 1).Select()
